# Auto-generated script to apply crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range('D2')
$cell.NumberFormat = '@'
$cell.Value = '67.788.12'
$cell.Style = 'Normal'

$cell = $ws.Range('D3')
$cell.NumberFormat = '@'
$cell.Value = '3.796.43'
$cell.Style = 'Normal'

$cell = $ws.Range('E3')
$cell.NumberFormat = '@'
$cell.Value = '  +0.20%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D4')
$cell.NumberFormat = '@'
$cell.Value = '0.999'
$cell.Style = 'Normal'

$cell = $ws.Range('E4')
$cell.NumberFormat = '@'
$cell.Value = '  -0.07%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '599.66'
$cell.Style = 'Normal'

$cell = $ws.Range('E5')
$cell.NumberFormat = '@'
$cell.Value = '  +0.71%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '165.29'
$cell.Style = 'Normal'

$cell = $ws.Range('E6')
$cell.NumberFormat = '@'
$cell.Value = '  -1.07%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E7')
$cell.NumberFormat = '@'
$cell.Value = '  +0.00%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E8')
$cell.NumberFormat = '@'
$cell.Value = '  -0.31%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E9')
$cell.NumberFormat = '@'
$cell.Value = '  -0.25%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D10')
$cell.NumberFormat = '@'
$cell.Value = '0.452'
$cell.Style = 'Normal'

$cell = $ws.Range('E10')
$cell.NumberFormat = '@'
$cell.Value = '  +0.92%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '6.48'
$cell.Style = 'Normal'

$cell = $ws.Range('E11')
$cell.NumberFormat = '@'
$cell.Value = '  +3.00%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E12')
$cell.NumberFormat = '@'
$cell.Value = '  -2.09%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D13')
$cell.NumberFormat = '@'
$cell.Value = '35.83'
$cell.Style = 'Normal'

$cell = $ws.Range('E13')
$cell.NumberFormat = '@'
$cell.Value = '  -0.73%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '4.436.13'
$cell.Style = 'Normal'

$cell = $ws.Range('E14')
$cell.NumberFormat = '@'
$cell.Value = '  +0.33%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '3.787.06'
$cell.Style = 'Normal'

$cell = $ws.Range('E15')
$cell.NumberFormat = '@'
$cell.Value = '  -0.09%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D16')
$cell.NumberFormat = '@'
$cell.Value = '67.834.28'
$cell.Style = 'Normal'

$cell = $ws.Range('E16')
$cell.NumberFormat = '@'
$cell.Value = '  +0.21%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D17')
$cell.NumberFormat = '@'
$cell.Value = '18.36'
$cell.Style = 'Normal'

$cell = $ws.Range('E17')
$cell.NumberFormat = '@'
$cell.Value = '  -0.82%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E18')
$cell.NumberFormat = '@'
$cell.Value = '  +1.84%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D19')
$cell.NumberFormat = '@'
$cell.Value = '7.06'
$cell.Style = 'Normal'

$cell = $ws.Range('E19')
$cell.NumberFormat = '@'
$cell.Value = '  +0.69%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '463.13'
$cell.Style = 'Normal'

$cell = $ws.Range('E20')
$cell.NumberFormat = '@'
$cell.Value = '  +0.83%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '9.81'
$cell.Style = 'Normal'

$cell = $ws.Range('E21')
$cell.NumberFormat = '@'
$cell.Value = '  -2.02%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '0.700'
$cell.Style = 'Normal'

$cell = $ws.Range('E22')
$cell.NumberFormat = '@'
$cell.Value = '  +0.34%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D23')
$cell.NumberFormat = '@'
$cell.Value = '0.0000146'
$cell.Style = 'Normal'

$cell = $ws.Range('E23')
$cell.NumberFormat = '@'
$cell.Value = '  -6.23%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D24')
$cell.NumberFormat = '@'
$cell.Value = '82.80'
$cell.Style = 'Normal'

$cell = $ws.Range('E24')
$cell.NumberFormat = '@'
$cell.Value = '  -0.57%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E25')
$cell.NumberFormat = '@'
$cell.Value = '  +0.20%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '2.10'
$cell.Style = 'Normal'

$cell = $ws.Range('E26')
$cell.NumberFormat = '@'
$cell.Value = '  -0.49%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '10.03'
$cell.Style = 'Normal'

$cell = $ws.Range('E27')
$cell.NumberFormat = '@'
$cell.Value = '  +0.05%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E28')
$cell.NumberFormat = '@'
$cell.Value = '  -0.17%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '3.946.42'
$cell.Style = 'Normal'

$cell = $ws.Range('E29')
$cell.NumberFormat = '@'
$cell.Value = '  +0.40%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E30')
$cell.NumberFormat = '@'
$cell.Value = '  -2.84%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '7.44'
$cell.Style = 'Normal'

$cell = $ws.Range('E31')
$cell.NumberFormat = '@'
$cell.Value = '  +2.85%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E32')
$cell.NumberFormat = '@'
$cell.Value = '  -1.43%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '29.19'
$cell.Style = 'Normal'

$cell = $ws.Range('E33')
$cell.NumberFormat = '@'
$cell.Value = '  -1.61%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D34')
$cell.NumberFormat = '@'
$cell.Value = '0.997'
$cell.Style = 'Normal'

$cell = $ws.Range('E34')
$cell.NumberFormat = '@'
$cell.Value = '  -0.26%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value = '9.03'
$cell.Style = 'Normal'

$cell = $ws.Range('E35')
$cell.NumberFormat = '@'
$cell.Value = '  -0.61%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '0.0994'
$cell.Style = 'Normal'

$cell = $ws.Range('E36')
$cell.NumberFormat = '@'
$cell.Value = '  -0.61%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E37')
$cell.NumberFormat = '@'
$cell.Value = '  +1.02%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D38')
$cell.NumberFormat = '@'
$cell.Value = '3.23'
$cell.Style = 'Normal'

$cell = $ws.Range('E38')
$cell.NumberFormat = '@'
$cell.Value = '  -4.66%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '5.77'
$cell.Style = 'Normal'

$cell = $ws.Range('E39')
$cell.NumberFormat = '@'
$cell.Value = '  -0.08%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D40')
$cell.NumberFormat = '@'
$cell.Value = '0.986'
$cell.Style = 'Normal'

$cell = $ws.Range('E40')
$cell.NumberFormat = '@'
$cell.Value = '  -0.67%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '1.00'
$cell.Style = 'Normal'

$cell = $ws.Range('E41')
$cell.NumberFormat = '@'
$cell.Value = '  +0.01%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E42')
$cell.NumberFormat = '@'
$cell.Value = '  +0.03%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D43')
$cell.NumberFormat = '@'
$cell.Value = '44.71'
$cell.Style = 'Normal'

$cell = $ws.Range('E43')
$cell.NumberFormat = '@'
$cell.Value = '  +0.27%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '47.59'
$cell.Style = 'Normal'

$cell = $ws.Range('E44')
$cell.NumberFormat = '@'
$cell.Value = '  -0.90%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E45')
$cell.NumberFormat = '@'
$cell.Value = '  +0.09%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D46')
$cell.NumberFormat = '@'
$cell.Value = '151.26'
$cell.Style = 'Normal'

$cell = $ws.Range('E46')
$cell.NumberFormat = '@'
$cell.Value = '  +1.17%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E47')
$cell.NumberFormat = '@'
$cell.Value = '  +8.82%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '8.35'
$cell.Style = 'Normal'

$cell = $ws.Range('E48')
$cell.NumberFormat = '@'
$cell.Value = '  +0.70%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '398.87'
$cell.Style = 'Normal'

$cell = $ws.Range('E49')
$cell.NumberFormat = '@'
$cell.Value = '  +1.02%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '27.36'
$cell.Style = 'Normal'

$cell = $ws.Range('E50')
$cell.NumberFormat = '@'
$cell.Value = '  +2.81%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '1.85'
$cell.Style = 'Normal'

$cell = $ws.Range('E51')
$cell.NumberFormat = '@'
$cell.Value = '  +1.68%  '
$cell.Style = 'Normal'

Write-Host "Applied 86 cell updates"